# redis-showcase.xlsx edit:
#  - add two new "#system" sheet columns: "sms" (before the old "sound"
#    column) and "ws.async" (before the old "xml" column), shifting the
#    columns in between to the right
#  - extend the "target" list (column A) with the two new category names
#  - update/add the workbook-level defined names so they keep pointing at
#    the right data after the shift

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("#system")

# ---------------------------------------------------------------------
# 1) Insert the two new columns. EntireColumn.Insert() shifts everything
#    at/after the target column one slot to the right, carrying existing
#    values (and their shared-string backing) along with it.
# ---------------------------------------------------------------------

# New column for "sms" goes in before the old column Q ("sound"); this
# pushes sound/ssh/step/web/webalert/webcookie/ws/xml from Q..X to R..Y.
$ws.Range("Q1").EntireColumn.Insert()

# New column for "ws.async" goes in before the (now shifted) "xml"
# column, which sits at Y after the previous insert; this pushes xml
# from Y to Z.
$ws.Range("Y1").EntireColumn.Insert()

# ---------------------------------------------------------------------
# 2) Populate the new "sms" column (Q) and its data.
# ---------------------------------------------------------------------
$ws.Range("Q1").Value = "sms"
$ws.Range("Q2").Value = "sendText(phones,text)"

# ---------------------------------------------------------------------
# 3) Populate the new "ws.async" column (Y) and its data.
# ---------------------------------------------------------------------
$ws.Range("Y1").Value = "ws.async"
$ws.Range("Y2").Value = "download(url,queryString,saveTo)"
$ws.Range("Y3").Value = "get(url,queryString,output)"
$ws.Range("Y4").Value = "head(url,output)"
$ws.Range("Y5").Value = "patch(url,body,output)"
$ws.Range("Y6").Value = "post(url,body,output)"
$ws.Range("Y7").Value = "put(url,body,output)"

# ---------------------------------------------------------------------
# 4) Column A ("target") lists the name of every category/defined name.
#    Insert "sms" and "ws.async" into that list in their proper spot,
#    shifting the trailing entries down. (Range.Insert on this engine
#    shifts the whole row, so we do the shift with plain value writes
#    instead, from the bottom up to avoid clobbering data.)
# ---------------------------------------------------------------------
$ws.Range("A26").Value = "xml"
$ws.Range("A25").Value = "ws.async"
$ws.Range("A24").Value = "ws"
$ws.Range("A23").Value = "webcookie"
$ws.Range("A22").Value = "webalert"
$ws.Range("A21").Value = "web"
$ws.Range("A20").Value = "step"
$ws.Range("A19").Value = "ssh"
$ws.Range("A18").Value = "sound"
$ws.Range("A17").Value = "sms"

# ---------------------------------------------------------------------
# 5) Fix up the defined names so they refer to the post-shift ranges,
#    and register the two new ones.
# ---------------------------------------------------------------------
$wb.Names.Item("sound").RefersTo     = "='#system'!`$R`$2:`$R`$5"
$wb.Names.Item("ssh").RefersTo       = "='#system'!`$S`$2:`$S`$9"
$wb.Names.Item("step").RefersTo      = "='#system'!`$T`$2:`$T`$4"
$wb.Names.Item("target").RefersTo    = "='#system'!`$A`$2:`$A`$26"
$wb.Names.Item("web").RefersTo       = "='#system'!`$U`$2:`$U`$108"
$wb.Names.Item("webalert").RefersTo  = "='#system'!`$V`$2:`$V`$6"
$wb.Names.Item("webcookie").RefersTo = "='#system'!`$W`$2:`$W`$8"
$wb.Names.Item("ws").RefersTo        = "='#system'!`$X`$2:`$X`$16"
$wb.Names.Item("xml").RefersTo       = "='#system'!`$Z`$2:`$Z`$11"

$wb.Names.Add("sms", "='#system'!`$Q`$2:`$Q`$2")
$wb.Names.Add("ws.async", "='#system'!`$Y`$2:`$Y`$7")
